$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Insert a new paragraph explaining the PlayerTurn/Dictionary
# design decision, right after the paragraph ending "...imitated with the
# text-based version of the game now." and before the "Game" heading.
# ---------------------------------------------------------------------------
$gameplayPara = $d.Paragraphs.Item(12)
$gameplayPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(13)

$playerTurnXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r><w:t xml:space="preserve">The </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>PlayerTurn</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> returning a Dictionary with two </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>ArrayLists</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> was made due to needing a compact way to return two separate </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>ArrayLists</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> from a single function. Since Tuples do not exist in Java, this was the solution that was settled upon.</w:t></w:r>' +
  '</w:p>'

$newPara.Range.InsertXML($playerTurnXml)

# ---------------------------------------------------------------------------
# Change 2: Mark "hard-coding" with gramStart/gramEnd proofErr tags inside
# the Words.txt explanation paragraph (splitting the run it lives in).
#
# This paragraph is the very last paragraph in the document body, and
# InsertXML on the range of the final paragraph does not replace it in
# place (it appends a sibling paragraph instead, leaving the original
# untouched). To work around this, temporarily add an empty paragraph
# after it so it is no longer "last", perform the replacement, then
# remove the now-unneeded trailing empty paragraph.
# ---------------------------------------------------------------------------
$wordsTxtPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$wordsTxtPara.Range.InsertParagraphAfter()
$wordsTxtPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)

$wordsTxtXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r><w:t xml:space="preserve">This text file contains all the valid words </w:t></w:r>' +
  '<w:r><w:t>allowed in the Scrabble game. It is loaded into the Board class' + [char]0x2019 + 's code at the start of the game as a Hash</w:t></w:r>' +
  '<w:r><w:t>Set</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">This text file is used due to the simplicity of editing it rather than </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>hard-coding</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve"> in adding every single one of the 10,000 possible words </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">manually into the HashSet. </w:t></w:r>' +
  '</w:p>'

$wordsTxtPara.Range.InsertXML($wordsTxtXml)

$trailingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailingPara.Range.Delete()

Write-Host "Edits applied"
